$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.001628901849080777
$ws.Range("J2").Value = 0.001628901849080777
$ws.Range("M2").Value = 19.21315233333334
$ws.Range("N2").Value = 57.63945700000001
$ws.Range("O2").Value = 0.04451179209991234
$ws.Range("P2").Value = 0.04451179209991233
$ws.Range("Q2").Value = 1.799228459023223
$ws.Range("R2").Value = 16.193056131209
$ws.Range("S2").Value = 0.00007250534045744633
$ws.Range("T2").Value = 0.00007250534045744632
$ws.Range("I3").Value = 0.001628901849080777
$ws.Range("J3").Value = 0.001628901849080777
$ws.Range("O3").Value = 0.2141755495962477
$ws.Range("P3").Value = 0.2141755495962477
$ws.Range("S3").Value = 0.0003488709487652196
$ws.Range("T3").Value = 0.0003488709487652195
$ws.Range("I4").Value = 0.001628901849080777
$ws.Range("J4").Value = 0.001628901849080777
$ws.Range("M4").Value = 166.8580016666666
$ws.Range("N4").Value = 500.5740049999999
$ws.Range("O4").Value = 0.3865658561145097
$ws.Range("P4").Value = 0.3865658561145097
$ws.Range("Q4").Value = 15.62552880474277
$ws.Range("R4").Value = 140.629759242685
$ws.Range("S4").Value = 0.0006296778378164185
$ws.Range("T4").Value = 0.0006296778378164184
$ws.Range("I5").Value = 0.001628901849080777
$ws.Range("J5").Value = 0.001628901849080777
$ws.Range("M5").Value = 41.09915599999999
$ws.Range("N5").Value = 123.297468
$ws.Range("O5").Value = 0.09521587377309249
$ws.Range("P5").Value = 0.09521587377309249
$ws.Range("Q5").Value = 3.848757863057333
$ws.Range("R5").Value = 34.63882076751599
$ws.Range("S5").Value = 0.0001550973128508322
$ws.Range("T5").Value = 0.0001550973128508322
$ws.Range("I6").Value = 0.001628901849080777
$ws.Range("J6").Value = 0.001628901849080777
$ws.Range("M6").Value = 112.0244103333333
$ws.Range("N6").Value = 336.073231
$ws.Range("O6").Value = 0.2595309284162377
$ws.Range("P6").Value = 0.2595309284162377
$ws.Range("Q6").Value = 10.49060058860522
$ws.Range("R6").Value = 94.41540529744701
$ws.Range("S6").Value = 0.0004227504091908604
$ws.Range("T6").Value = 0.0004227504091908603
$ws.Range("I7").Value = 0.1785947081647151
$ws.Range("J7").Value = 0.178594708164715
$ws.Range("M7").Value = 19.21315233333334
$ws.Range("N7").Value = 57.63945700000001
$ws.Range("O7").Value = 0.04451179209991234
$ws.Range("P7").Value = 0.04451179209991233
$ws.Range("Q7").Value = 197.2695173390816
$ws.Range("R7").Value = 1775.425656051734
$ws.Range("S7").Value = 0.007949570519972313
$ws.Range("T7").Value = 0.007949570519972311
$ws.Range("I8").Value = 0.1785947081647151
$ws.Range("J8").Value = 0.178594708164715
$ws.Range("O8").Value = 0.2141755495962477
$ws.Range("P8").Value = 0.2141755495962477
$ws.Range("S8").Value = 0.03825061977615932
$ws.Range("T8").Value = 0.03825061977615931
$ws.Range("I9").Value = 0.1785947081647151
$ws.Range("J9").Value = 0.178594708164715
$ws.Range("M9").Value = 166.8580016666666
$ws.Range("N9").Value = 500.5740049999999
$ws.Range("O9").Value = 0.3865658561145097
$ws.Range("P9").Value = 0.3865658561145097
$ws.Range("Q9").Value = 1713.201294711034
$ws.Range("R9").Value = 15418.81165239931
$ws.Range("S9").Value = 0.06903861625921409
$ws.Range("T9").Value = 0.06903861625921409
$ws.Range("I10").Value = 0.1785947081647151
$ws.Range("J10").Value = 0.178594708164715
$ws.Range("M10").Value = 41.09915599999999
$ws.Range("N10").Value = 123.297468
$ws.Range("O10").Value = 0.09521587377309249
$ws.Range("P10").Value = 0.09521587377309249
$ws.Range("Q10").Value = 421.9823236969573
$ws.Range("R10").Value = 3797.840913272615
$ws.Range("S10").Value = 0.0170050511891538
$ws.Range("T10").Value = 0.0170050511891538
$ws.Range("I11").Value = 0.1785947081647151
$ws.Range("J11").Value = 0.178594708164715
$ws.Range("M11").Value = 112.0244103333333
$ws.Range("N11").Value = 336.073231
$ws.Range("O11").Value = 0.2595309284162377
$ws.Range("P11").Value = 0.2595309284162377
$ws.Range("Q11").Value = 1150.201745827614
$ws.Range("R11").Value = 10351.81571244852
$ws.Range("S11").Value = 0.04635085042021553
$ws.Range("T11").Value = 0.04635085042021552
$ws.Range("G12").Value = 23.67539566666666
$ws.Range("H12").Value = 71.02618699999999
$ws.Range("I12").Value = 0.4118171950916292
$ws.Range("J12").Value = 0.4118171950916292
$ws.Range("M12").Value = 19.21315233333334
$ws.Range("N12").Value = 57.63945700000001
$ws.Range("O12").Value = 0.04451179209991234
$ws.Range("P12").Value = 0.04451179209991233
$ws.Range("Q12").Value = 454.8789834956066
$ws.Range("R12").Value = 4093.910851460459
$ws.Range("S12").Value = 0.01833072137108764
$ws.Range("T12").Value = 0.01833072137108763
$ws.Range("G13").Value = 23.67539566666666
$ws.Range("H13").Value = 71.02618699999999
$ws.Range("I13").Value = 0.4118171950916292
$ws.Range("J13").Value = 0.4118171950916292
$ws.Range("O13").Value = 0.2141755495962477
$ws.Range("P13").Value = 0.2141755495962477
$ws.Range("Q13").Value = 2188.722396781367
$ws.Range("R13").Value = 19698.5015710323
$ws.Range("S13").Value = 0.08820117409193487
$ws.Range("T13").Value = 0.08820117409193484
$ws.Range("G14").Value = 23.67539566666666
$ws.Range("H14").Value = 71.02618699999999
$ws.Range("I14").Value = 0.4118171950916292
$ws.Range("J14").Value = 0.4118171950916292
$ws.Range("M14").Value = 166.8580016666666
$ws.Range("N14").Value = 500.5740049999999
$ws.Range("O14").Value = 0.3865658561145097
$ws.Range("P14").Value = 0.3865658561145097
$ws.Range("Q14").Value = 3950.429209607658
$ws.Range("R14").Value = 35553.86288646893
$ws.Range("S14").Value = 0.1591944665832717
$ws.Range("T14").Value = 0.1591944665832717
$ws.Range("G15").Value = 23.67539566666666
$ws.Range("H15").Value = 71.02618699999999
$ws.Range("I15").Value = 0.4118171950916292
$ws.Range("J15").Value = 0.4118171950916292
$ws.Range("M15").Value = 41.09915599999999
$ws.Range("N15").Value = 123.297468
$ws.Range("O15").Value = 0.09521587377309249
$ws.Range("P15").Value = 0.09521587377309249
$ws.Range("Q15").Value = 973.038779866057
$ws.Range("R15").Value = 8757.349018794514
$ws.Range("S15").Value = 0.03921153406543357
$ws.Range("T15").Value = 0.03921153406543357
$ws.Range("G16").Value = 23.67539566666666
$ws.Range("H16").Value = 71.02618699999999
$ws.Range("I16").Value = 0.4118171950916292
$ws.Range("J16").Value = 0.4118171950916292
$ws.Range("M16").Value = 112.0244103333333
$ws.Range("N16").Value = 336.073231
$ws.Range("O16").Value = 0.2595309284162377
$ws.Range("P16").Value = 0.2595309284162377
$ws.Range("Q16").Value = 2652.222238966688
$ws.Range("R16").Value = 23870.0001507002
$ws.Range("S16").Value = 0.1068792989799014
$ws.Range("T16").Value = 0.1068792989799014
$ws.Range("G17").Value = 0.3314846666666666
$ws.Range("H17").Value = 0.9944539999999999
$ws.Range("I17").Value = 0.005765947381177186
$ws.Range("J17").Value = 0.005765947381177185
$ws.Range("M17").Value = 19.21315233333334
$ws.Range("N17").Value = 57.63945700000001
$ws.Range("O17").Value = 0.04451179209991234
$ws.Range("P17").Value = 0.04451179209991233
$ws.Range("Q17").Value = 6.368865396830889
$ws.Range("R17").Value = 57.31978857147801
$ws.Range("S17").Value = 0.0002566526510899929
$ws.Range("T17").Value = 0.0002566526510899928
$ws.Range("G18").Value = 0.3314846666666666
$ws.Range("H18").Value = 0.9944539999999999
$ws.Range("I18").Value = 0.005765947381177186
$ws.Range("J18").Value = 0.005765947381177185
$ws.Range("O18").Value = 0.2141755495962477
$ws.Range("P18").Value = 0.2141755495962477
$ws.Range("Q18").Value = 30.64480629332978
$ws.Range("R18").Value = 275.803256639968
$ws.Range("S18").Value = 0.001234924949306669
$ws.Range("T18").Value = 0.001234924949306669
$ws.Range("G19").Value = 0.3314846666666666
$ws.Range("H19").Value = 0.9944539999999999
$ws.Range("I19").Value = 0.005765947381177186
$ws.Range("J19").Value = 0.005765947381177185
$ws.Range("M19").Value = 166.8580016666666
$ws.Range("N19").Value = 500.5740049999999
$ws.Range("O19").Value = 0.3865658561145097
$ws.Range("P19").Value = 0.3865658561145097
$ws.Range("Q19").Value = 55.3108690631411
$ws.Range("R19").Value = 497.7978215682699
$ws.Range("S19").Value = 0.002228918385715974
$ws.Range("T19").Value = 0.002228918385715974
$ws.Range("G20").Value = 0.3314846666666666
$ws.Range("H20").Value = 0.9944539999999999
$ws.Range("I20").Value = 0.005765947381177186
$ws.Range("J20").Value = 0.005765947381177185
$ws.Range("M20").Value = 41.09915599999999
$ws.Range("N20").Value = 123.297468
$ws.Range("O20").Value = 0.09521587377309249
$ws.Range("P20").Value = 0.09521587377309249
$ws.Range("Q20").Value = 13.62374002694133
$ws.Range("R20").Value = 122.613660242472
$ws.Range("S20").Value = 0.0005490097180284601
$ws.Range("T20").Value = 0.00054900971802846
$ws.Range("G21").Value = 0.3314846666666666
$ws.Range("H21").Value = 0.9944539999999999
$ws.Range("I21").Value = 0.005765947381177186
$ws.Range("J21").Value = 0.005765947381177185
$ws.Range("M21").Value = 112.0244103333333
$ws.Range("N21").Value = 336.073231
$ws.Range("O21").Value = 0.2595309284162377
$ws.Range("P21").Value = 0.2595309284162377
$ws.Range("Q21").Value = 37.13437431787489
$ws.Range("R21").Value = 334.209368860874
$ws.Range("S21").Value = 0.001496441677036089
$ws.Range("T21").Value = 0.001496441677036089
$ws.Range("G22").Value = 23.12211433333333
$ws.Range("H22").Value = 69.366343
$ws.Range("I22").Value = 0.4021932475133977
$ws.Range("J22").Value = 0.4021932475133977
$ws.Range("M22").Value = 19.21315233333334
$ws.Range("N22").Value = 57.63945700000001
$ws.Range("O22").Value = 0.04451179209991234
$ws.Range("P22").Value = 0.04451179209991233
$ws.Range("Q22").Value = 444.2487049550835
$ws.Range("R22").Value = 3998.238344595752
$ws.Range("S22").Value = 0.01790234221730495
$ws.Range("T22").Value = 0.01790234221730494
$ws.Range("G23").Value = 23.12211433333333
$ws.Range("H23").Value = 69.366343
$ws.Range("I23").Value = 0.4021932475133977
$ws.Range("J23").Value = 0.4021932475133977
$ws.Range("O23").Value = 0.2141755495962477
$ws.Range("P23").Value = 0.2141755495962477
$ws.Range("Q23").Value = 2137.573125063273
$ws.Range("R23").Value = 19238.15812556945
$ws.Range("S23").Value = 0.08613995983008166
$ws.Range("T23").Value = 0.08613995983008164
$ws.Range("G24").Value = 23.12211433333333
$ws.Range("H24").Value = 69.366343
$ws.Range("I24").Value = 0.4021932475133977
$ws.Range("J24").Value = 0.4021932475133977
$ws.Range("M24").Value = 166.8580016666666
$ws.Range("N24").Value = 500.5740049999999
$ws.Range("O24").Value = 0.3865658561145097
$ws.Range("P24").Value = 0.3865658561145097
$ws.Range("Q24").Value = 3858.10979196819
$ws.Range("R24").Value = 34722.98812771371
$ws.Range("S24").Value = 0.1554741770484915
$ws.Range("T24").Value = 0.1554741770484915
$ws.Range("G25").Value = 23.12211433333333
$ws.Range("H25").Value = 69.366343
$ws.Range("I25").Value = 0.4021932475133977
$ws.Range("J25").Value = 0.4021932475133977
$ws.Range("M25").Value = 41.09915599999999
$ws.Range("N25").Value = 123.297468
$ws.Range("O25").Value = 0.09521587377309249
$ws.Range("P25").Value = 0.09521587377309249
$ws.Range("Q25").Value = 950.2993840355025
$ws.Range("R25").Value = 8552.694456319523
$ws.Range("S25").Value = 0.03829518148762583
$ws.Range("T25").Value = 0.03829518148762582
$ws.Range("G26").Value = 23.12211433333333
$ws.Range("H26").Value = 69.366343
$ws.Range("I26").Value = 0.4021932475133977
$ws.Range("J26").Value = 0.4021932475133977
$ws.Range("M26").Value = 112.0244103333333
$ws.Range("N26").Value = 336.073231
$ws.Range("O26").Value = 0.2595309284162377
$ws.Range("P26").Value = 0.2595309284162377
$ws.Range("Q26").Value = 2590.241223851581
$ws.Range("R26").Value = 23312.17101466424
$ws.Range("S26").Value = 0.1043815869298938
$ws.Range("T26").Value = 0.1043815869298938
